# Update TPM-derived values in the LR-pairs sheet (Fndc5-Itgav) to reflect
# the newly computed statistics from the updated TPM script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "G2" = 0.15188
    "H2" = 0.45564
    "I2" = 0.229582042173683
    "J2" = 0.229582042173683
    "M2" = 3.759736666666667
    "N2" = 11.27921
    "O2" = 0.0683751702595819
    "P2" = 0.06837517025958188
    "Q2" = 0.5710288049333333
    "R2" = 5.1392592444
    "S2" = 0.01569771122216809
    "T2" = 0.01569771122216808

    "G3" = 0.15188
    "H3" = 0.45564
    "I3" = 0.229582042173683
    "J3" = 0.229582042173683
    "O3" = 0.6514180024294648
    "P3" = 0.6514180024294647
    "Q3" = 5.440256192813333
    "R3" = 48.96230573532
    "S3" = 0.1495538753064577
    "T3" = 0.1495538753064577

    "G4" = 0.15188
    "H4" = 0.45564
    "I4" = 0.229582042173683
    "J4" = 0.229582042173683
    "O4" = 0.2802068273109533
    "P4" = 0.2802068273109533
    "Q4" = 2.340120969733333
    "R4" = 21.0610887276
    "S4" = 0.0643304556450572
    "T4" = 0.06433045564505718

    "I5" = 0.7704179578263169
    "J5" = 0.7704179578263169
    "M5" = 3.759736666666667
    "N5" = 11.27921
    "O5" = 0.0683751702595819
    "P5" = 0.06837517025958188
    "Q5" = 1.9162249869
    "R5" = 17.2460248821
    "S5" = 0.0526774590374138
    "T5" = 0.0526774590374138

    "I6" = 0.7704179578263169
    "J6" = 0.7704179578263169
    "O6" = 0.6514180024294648
    "P6" = 0.6514180024294647
    "S6" = 0.501864127123007
    "T6" = 0.501864127123007

    "I7" = 0.7704179578263169
    "J7" = 0.7704179578263169
    "O7" = 0.2802068273109533
    "P7" = 0.2802068273109533
    "S7" = 0.2158763716658961
    "T7" = 0.2158763716658961
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
